$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data column (BQ, "14-sep") into the
# new column (BR) so the new "15-sep" column matches the existing look
# (text-style header, centered integer body cells) before filling in values.
$ws.Range("BQ1:BQ11").Copy()
$ws.Range("BR1:BR11").PasteSpecial(-4122)

# New column header: 15-sep
$ws.Range("BR1").Value = "15-sep"

# New column values (one per product row)
$ws.Range("BR2").Value = 15
$ws.Range("BR3").Value = 13
$ws.Range("BR4").Value = 12
$ws.Range("BR5").Value = 11
$ws.Range("BR6").Value = 12
$ws.Range("BR7").Value = 17
$ws.Range("BR8").Value = 19
$ws.Range("BR9").Value = 10
$ws.Range("BR10").Value = 12
$ws.Range("BR11").Value = 10

# Match the author's final selection/active cell after entering the new column
$ws.Range("BR12").Select()
